$d = $word.ActiveDocument

# Locate the paragraphs that need a new bold heading inserted right before them,
# by searching for their distinctive text (robust to any pre-existing paragraph
# ordering assumptions).
$find1 = $d.Content
$ok1 = $find1.Find.Execute("El cliente solicita", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$idxSolicita = $find1.Paragraphs.First.Index

$find2 = $d.Content
$ok2 = $find2.Find.Execute("El cliente carece", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$idxCarece = $find2.Paragraphs.First.Index

# --- Insert "Introducción" heading before the "El cliente solicita..." paragraph ---
$ins1 = $d.Paragraphs.Item($idxSolicita).Range.InsertParagraphBefore()

$newPara1 = $d.Paragraphs.Item($idxSolicita)
$newPara1.Range.Text = "Introducción"
$newRun1 = $d.Range($newPara1.Range.Start, $newPara1.Range.Start + 12)
$newRun1.Font.Bold = 1

# --- Insert "Enunciado del problema" heading before the "El cliente carece..." paragraph ---
# (shift by 1 because the paragraph inserted above moved everything after it down by one)
$ins2 = $d.Paragraphs.Item($idxCarece + 1).Range.InsertParagraphBefore()

$newPara2 = $d.Paragraphs.Item($idxCarece + 1)
$newPara2.Range.Text = "Enunciado del problema"
$newRun2 = $d.Range($newPara2.Range.Start, $newPara2.Range.Start + 23)
$newRun2.Font.Bold = 1

Write-Output "done"
